# New crime data collected - weekly CompStat update (60th Precinct)
# Updates: volume/issue number, week-covering dates, and the weekly crime
# statistics table (rows 14-30) with newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header: "Volume 32   Number  48" -> "Volume 32   Number  49"
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  49"

# ---------------------------------------------------------------------
# Header: "Report Covering the Week  11/24/2025  Through  11/30/2025"
#      -> "Report Covering the Week  12/1/2025  Through  12/7/2025"
# ---------------------------------------------------------------------
$ws.Range("C9").Value = "Report Covering the Week  12/1/2025  Through  12/7/2025"

# ---------------------------------------------------------------------
# Crime-complaint grid (rows 14-30, columns C..N)
# Row 14 (Murder) has no new data this week.
# ---------------------------------------------------------------------

# Row 15 - Rape
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("N15").Value = 11.111111111111

# Row 16 - Robbery
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 146
$ws.Range("J16").Value = 169
$ws.Range("K16").Value = -13.609467455621
$ws.Range("L16").Value = -8.176100628930
$ws.Range("M16").Value = -42.745098039215
$ws.Range("N16").Value = -85.472636815920

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 47.368421052631
$ws.Range("I17").Value = 409
$ws.Range("J17").Value = 339
$ws.Range("K17").Value = 20.648967551622
$ws.Range("L17").Value = 17.191977077363
$ws.Range("M17").Value = 102.475247524752
$ws.Range("N17").Value = -41.985815602836

# Row 18 - Burglary (C18 switches from "n/a" text to a real number)
$ws.Range("C18").Value = 4
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 110
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 10
$ws.Range("L18").Value = 11.111111111111
$ws.Range("M18").Value = -29.032258064516
$ws.Range("N18").Value = -88.159311087190

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 22
$ws.Range("H19").Value = 59.090909090909
$ws.Range("I19").Value = 397
$ws.Range("J19").Value = 334
$ws.Range("K19").Value = 18.862275449101
$ws.Range("L19").Value = -10.383747178329
$ws.Range("M19").Value = -24.236641221374
$ws.Range("N19").Value = -39.016897081413

# Row 20 - G.L.A.
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -80
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -35.714285714285
$ws.Range("I20").Value = 101
$ws.Range("J20").Value = 112
$ws.Range("K20").Value = -9.821428571428
$ws.Range("L20").Value = -4.716981132075
$ws.Range("M20").Value = -6.481481481481
$ws.Range("N20").Value = -90.876242095754

# Row 21 - TOTAL
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 31.578947368421
$ws.Range("F21").Value = 99
$ws.Range("G21").Value = 80
$ws.Range("H21").Value = 23.75
$ws.Range("I21").Value = 1207
$ws.Range("J21").Value = 1074
$ws.Range("K21").Value = 12.383612662942
$ws.Range("L21").Value = 3.074295473953
$ws.Range("M21").Value = -4.810725552050
$ws.Range("N21").Value = -72.894677745340

# Row 22 - Transit (C22 switches from "n/a" text to a real number)
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("I22").Value = 45
$ws.Range("J22").Value = 28
$ws.Range("K22").Value = 60.714285714285
$ws.Range("L22").Value = 114.285714285714
$ws.Range("M22").Value = 55.172413793103

# Row 23 - Housing
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 33.333333333333
$ws.Range("I23").Value = 144
$ws.Range("J23").Value = 131
$ws.Range("K23").Value = 9.923664122137
$ws.Range("L23").Value = -7.096774193548
$ws.Range("M23").Value = 25.217391304347

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 91.666666666666
$ws.Range("F24").Value = 106
$ws.Range("G24").Value = 65
$ws.Range("H24").Value = 63.076923076923
$ws.Range("I24").Value = 970
$ws.Range("J24").Value = 934
$ws.Range("K24").Value = 3.854389721627
$ws.Range("L24").Value = 4.301075268817
$ws.Range("M24").Value = -14.235190097259

# Row 25 - Retail Theft
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 250
$ws.Range("F25").Value = 21
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = -12.5
$ws.Range("I25").Value = 229
$ws.Range("J25").Value = 264
$ws.Range("K25").Value = -13.257575757575
$ws.Range("L25").Value = -17.625899280575

# Row 26 - Misd. Assault
$ws.Range("D26").Value = 17
$ws.Range("E26").Value = -29.411764705882
$ws.Range("F26").Value = 49
$ws.Range("G26").Value = 55
$ws.Range("H26").Value = -10.909090909090
$ws.Range("I26").Value = 729
$ws.Range("J26").Value = 743
$ws.Range("K26").Value = -1.884253028263
$ws.Range("L26").Value = 26.782608695652
$ws.Range("M26").Value = 46.975806451612

# Row 27 - UCR Rape* (C27 switches from "n/a" text to a real number)
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("I27").Value = 45
$ws.Range("K27").Value = 73.076923076923
$ws.Range("L27").Value = 181.25

# Row 28 - Other Sex Crimes (C28/D28 switch from "n/a" text to numbers,
# E28 switches from "n/a" text to a percentage value)
$ws.Range("C28").Value = 2
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = 100
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 69
$ws.Range("J28").Value = 45
$ws.Range("K28").Value = 53.333333333333
$ws.Range("L28").Value = 38

# Row 29 - Shooting Vic.
$ws.Range("L29").Value = -9.090909090909

# Row 30 - Shooting Inc.
$ws.Range("L30").Value = 28.571428571428
